# Rename the worksheet to match its actual contents (rainfall data).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "rainfall_data"

# Give the value column a descriptive header instead of the raw sensor id.
$ws.Range("B1").Value = "rainfall_value"

# The header text got longer, so the column needs to be widened to fit it
# (mirrors Excel's "best fit" auto-resize after the content changed).
$ws.Columns.Item(2).ColumnWidth = 12.45

# Move the active selection before saving.
$ws.Range("F8").Select()
